# Auto-generated edit script applying scheduled market-data refresh to Sheets/Coeurl_Profits.xlsx
# Updates currentAveragePrice/LevePrice/LeveProfit columns (H-N) for affected leve rows across sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3525
$ws.Range("I18").Value = 1366.6666
$ws.Range("K18").Value = 1366.6666
$ws.Range("M18").Value = -1082.6666
$ws.Range("H70").Value = 3954.524
$ws.Range("I70").Value = 4736.077
$ws.Range("K70").Value = 14208.231
$ws.Range("M70").Value = -13938.231
$ws.Range("H73").Value = 3954.524
$ws.Range("I73").Value = 4736.077
$ws.Range("K73").Value = 14208.231
$ws.Range("M73").Value = -13272.231
$ws.Range("H137").Value = 2018.5686
$ws.Range("I137").Value = 1961.6154
$ws.Range("J137").Value = 2203.6667
$ws.Range("K137").Value = 5884.8462
$ws.Range("L137").Value = 6611.000100000001
$ws.Range("M137").Value = -3334.8462
$ws.Range("N137").Value = -11711.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2385.5386
$ws.Range("I132").Value = 2127.5293
$ws.Range("K132").Value = 6382.5879
$ws.Range("M132").Value = -3852.5879

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3249.2856
$ws.Range("I20").Value = 3090.077
$ws.Range("J20").Value = 3508
$ws.Range("K20").Value = 3090.077
$ws.Range("L20").Value = 3508
$ws.Range("M20").Value = -2843.077
$ws.Range("N20").Value = -4002
$ws.Range("H94").Value = 1160.7
$ws.Range("I94").Value = 1333.9333
$ws.Range("K94").Value = 1333.9333
$ws.Range("M94").Value = -882.9332999999999
$ws.Range("H133").Value = 21542.334
$ws.Range("I133").Value = 21542.334
$ws.Range("K133").Value = 21542.334
$ws.Range("M133").Value = -16482.334
$ws.Range("H134").Value = 3374.5715
$ws.Range("I134").Value = 3374.5715
$ws.Range("K134").Value = 10123.7145
$ws.Range("M134").Value = -7588.7145

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 741
$ws.Range("I22").Value = 168.83333
$ws.Range("K22").Value = 168.83333
$ws.Range("M22").Value = 181.16667
$ws.Range("H41").Value = 15394.5
$ws.Range("J41").Value = 15789
$ws.Range("L41").Value = 15789
$ws.Range("N41").Value = -16645
$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -9376
$ws.Range("N62").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -46880
$ws.Range("N65").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H69").Value = 14795.5
$ws.Range("J69").Value = 17500
$ws.Range("L69").Value = 17500
$ws.Range("N69").Value = -18998
$ws.Range("H72").Value = 14795.5
$ws.Range("J72").Value = 17500
$ws.Range("L72").Value = 17500
$ws.Range("N72").Value = -59988
$ws.Range("H99").Value = 3068.182
$ws.Range("I99").Value = 2916.6667
$ws.Range("K99").Value = 2916.6667
$ws.Range("M99").Value = -1418.6667
$ws.Range("H126").Value = 3068.182
$ws.Range("I126").Value = 2916.6667
$ws.Range("K126").Value = 8750.000100000001
$ws.Range("M126").Value = -6280.000100000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 202.5
$ws.Range("I20").Value = 200
$ws.Range("J20").Value = 205
$ws.Range("K20").Value = 600
$ws.Range("L20").Value = 615
$ws.Range("M20").Value = -373
$ws.Range("N20").Value = -1069
$ws.Range("H64").Value = 6382.625
$ws.Range("J64").Value = 7512
$ws.Range("L64").Value = 22536
$ws.Range("N64").Value = -23076
$ws.Range("H67").Value = 6382.625
$ws.Range("J67").Value = 7512
$ws.Range("L67").Value = 22536
$ws.Range("N67").Value = -24408
$ws.Range("H100").Value = 10000
$ws.Range("J100").Value = 10000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -31622
$ws.Range("H107").Value = 733.1667
$ws.Range("I107").Value = 366
$ws.Range("K107").Value = 1098
$ws.Range("M107").Value = 822
$ws.Range("H112").Value = 5488.25
$ws.Range("I112").Value = 977.5
$ws.Range("K112").Value = 2932.5
$ws.Range("M112").Value = -1824.5
$ws.Range("H113").Value = 1000
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 7077.5264
$ws.Range("I15").Value = 10249.125
$ws.Range("J15").Value = 4770.909
$ws.Range("K15").Value = 10249.125
$ws.Range("L15").Value = 4770.909
$ws.Range("M15").Value = -9961.125
$ws.Range("N15").Value = -5346.909
$ws.Range("H81").Value = 7077.5264
$ws.Range("I81").Value = 10249.125
$ws.Range("J81").Value = 4770.909
$ws.Range("K81").Value = 10249.125
$ws.Range("L81").Value = 4770.909
$ws.Range("M81").Value = -9251.125
$ws.Range("N81").Value = -6766.909
$ws.Range("H84").Value = 7077.5264
$ws.Range("I84").Value = 10249.125
$ws.Range("J84").Value = 4770.909
$ws.Range("K84").Value = 30747.375
$ws.Range("L84").Value = 14312.727
$ws.Range("M84").Value = -25755.375
$ws.Range("N84").Value = -24296.727
$ws.Range("H132").Value = 3065.8262
$ws.Range("I132").Value = 2550.7
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 7652.099999999999
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -5122.099999999999
$ws.Range("N132").Value = -24560
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4725.0835
$ws.Range("I7").Value = 3378.353
$ws.Range("J7").Value = 7995.7144
$ws.Range("K7").Value = 3378.353
$ws.Range("L7").Value = 7995.7144
$ws.Range("M7").Value = -3266.353
$ws.Range("N7").Value = -8219.714400000001
$ws.Range("H46").Value = 1183.1875
$ws.Range("I46").Value = 685.8570999999999
$ws.Range("J46").Value = 1570
$ws.Range("K46").Value = 685.8570999999999
$ws.Range("L46").Value = 1570
$ws.Range("M46").Value = -497.8570999999999
$ws.Range("N46").Value = -1946
$ws.Range("H68").Value = 2754.5715
$ws.Range("I68").Value = 2833.7368
$ws.Range("J68").Value = 2002.5
$ws.Range("K68").Value = 2833.7368
$ws.Range("L68").Value = 2002.5
$ws.Range("M68").Value = -2084.7368
$ws.Range("N68").Value = -3500.5
$ws.Range("H71").Value = 2754.5715
$ws.Range("I71").Value = 2833.7368
$ws.Range("J71").Value = 2002.5
$ws.Range("K71").Value = 14168.684
$ws.Range("L71").Value = 10012.5
$ws.Range("M71").Value = -10424.684
$ws.Range("N71").Value = -17500.5
$ws.Range("H82").Value = 8572.235000000001
$ws.Range("I82").Value = 10171.385
$ws.Range("J82").Value = 3375
$ws.Range("K82").Value = 10171.385
$ws.Range("L82").Value = 3375
$ws.Range("M82").Value = -9810.385
$ws.Range("N82").Value = -4097
$ws.Range("H85").Value = 8572.235000000001
$ws.Range("I85").Value = 10171.385
$ws.Range("J85").Value = 3375
$ws.Range("K85").Value = 10171.385
$ws.Range("L85").Value = 3375
$ws.Range("M85").Value = -8923.385
$ws.Range("N85").Value = -5871
$ws.Range("H126").Value = 4725.0835
$ws.Range("I126").Value = 3378.353
$ws.Range("J126").Value = 7995.7144
$ws.Range("K126").Value = 10135.059
$ws.Range("L126").Value = 23987.1432
$ws.Range("M126").Value = -7665.059000000001
$ws.Range("N126").Value = -28927.1432

